# SD1 Meeting Log.xlsx update
# Adds two new meeting-log rows (2/8/2017 and 2/15/2017), updates the
# "Last updated" banner in B1, and formats the two new "Actions Required"
# cells as rich text (bold / bold+italic headers, as in the existing log).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Update the "Last updated" banner (B1)
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Last updated: 02/15/2017 12:28 PM by Carlos"

# ---------------------------------------------------------------------
# Row 6 - meeting of 2/8/2017 (serial date 42774)
# ---------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 42774

$ws.Range("B6").Value = "Meet with Lei Wei, choose project, Assign research"
$ws.Range("B6").WrapText = $true

$ws.Range("C6").Value = "Carlos, Courtnie, Lucas, Patrick"

$ws.Range("D6").Value = "Chose Fire Alarm Project. See ""Actions Required"" for Research."
$ws.Range("D6").WrapText = $true

$ws.Range("F6").Value = "Meeting Weekly on Wednesdays from now on"
$ws.Range("F6").WrapText = $true

# E6 - rich text "Actions Required" for the 2/8 meeting
$c = $ws.Range("E6")
$c.Value = 'Please complete the following by 2/15: Carlos: Update Project Milestones; Courtnie: Update House of Quality; Lucas: goals and objectives;  Patrick: Project Block Diagram; Complete by 2/22: Carlos: Smoke Sensors Research & Docs; Courtnie: Battery Research & Docs; Lucas: Microcontroller Research & Docs; Patrick: Wireless Communication Research & Docs;'
$c.WrapText = $true
$c.Font.Bold = $true
$c.Font.Italic = $true
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Characters(1,37).Font.Bold = $true
$c.Characters(1,37).Font.Italic = $true
$c.Characters(38,1).Font.Bold = $true
$c.Characters(38,1).Font.Italic = $false
$c.Characters(39,1).Font.Bold = $false
$c.Characters(39,1).Font.Italic = $false
$c.Characters(40,7).Font.Bold = $true
$c.Characters(40,7).Font.Italic = $false
$c.Characters(47,28).Font.Bold = $false
$c.Characters(47,28).Font.Italic = $false
$c.Characters(75,10).Font.Bold = $true
$c.Characters(75,10).Font.Italic = $false
$c.Characters(85,25).Font.Bold = $false
$c.Characters(85,25).Font.Italic = $false
$c.Characters(110,6).Font.Bold = $true
$c.Characters(110,6).Font.Italic = $false
$c.Characters(116,24).Font.Bold = $false
$c.Characters(116,24).Font.Italic = $false
$c.Characters(140,8).Font.Bold = $true
$c.Characters(140,8).Font.Italic = $false
$c.Characters(148,24).Font.Bold = $false
$c.Characters(148,24).Font.Italic = $false
$c.Characters(172,16).Font.Bold = $true
$c.Characters(172,16).Font.Italic = $true
$c.Characters(188,2).Font.Bold = $false
$c.Characters(188,2).Font.Italic = $false
$c.Characters(190,6).Font.Bold = $true
$c.Characters(190,6).Font.Italic = $false
$c.Characters(196,33).Font.Bold = $false
$c.Characters(196,33).Font.Italic = $false
$c.Characters(229,8).Font.Bold = $true
$c.Characters(229,8).Font.Italic = $false
$c.Characters(237,27).Font.Bold = $false
$c.Characters(237,27).Font.Italic = $false
$c.Characters(264,6).Font.Bold = $true
$c.Characters(264,6).Font.Italic = $false
$c.Characters(270,34).Font.Bold = $false
$c.Characters(270,34).Font.Italic = $false
$c.Characters(304,7).Font.Bold = $true
$c.Characters(304,7).Font.Italic = $false
$c.Characters(311,41).Font.Bold = $false
$c.Characters(311,41).Font.Italic = $false

$ws.Rows.Item(6).RowHeight = 90

# ---------------------------------------------------------------------
# Row 7 - meeting of 2/15/2017 (serial date 42781)
# ---------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 42781

$ws.Range("B7").Value = "Discuss Updates for Initial Document/Proposal. Make sure everyone has Research Assigned, and Understands what is Required."
$ws.Range("B7").WrapText = $true

$ws.Range("C7").Value = "Carlos, Courtnie, Lucas, Patrick"

$ws.Range("D7").Value = "Assigned Research for everyone. Will submit ""Updates"" documented by 2/16."
$ws.Range("D7").WrapText = $true

$ws.Range("F7").Value = "Will look into sponsors later on in project. Add Software Flowchart to document whenever possible."
$ws.Range("F7").WrapText = $true

# E7 - rich text "Actions Required" for the 2/15 meeting
$c = $ws.Range("E7")
$c.Value = 'Please complete Research on at least one of the following by 2/22: Carlos: Smoke Sensors, Direction/Location Algorithm Research & Docs; Courtnie: Battery/Power Monitoring, Alarm System Components Research & Docs; Lucas: Microcontroller and Processor Research & Docs; Patrick: Wireless Communication Research and Processor & Docs;'
$c.WrapText = $true
$c.Font.Bold = $true
$c.Font.Italic = $true
$c.Font.Bold = $false
$c.Font.Italic = $false
$c.Characters(1,65).Font.Bold = $true
$c.Characters(1,65).Font.Italic = $true
$c.Characters(66,8).Font.Bold = $true
$c.Characters(66,8).Font.Italic = $false
$c.Characters(74,63).Font.Bold = $false
$c.Characters(74,63).Font.Italic = $false
$c.Characters(137,8).Font.Bold = $true
$c.Characters(137,8).Font.Italic = $false
$c.Characters(145,69).Font.Bold = $false
$c.Characters(145,69).Font.Italic = $false
$c.Characters(214,6).Font.Bold = $true
$c.Characters(214,6).Font.Italic = $false
$c.Characters(220,48).Font.Bold = $false
$c.Characters(220,48).Font.Italic = $false
$c.Characters(268,7).Font.Bold = $true
$c.Characters(268,7).Font.Italic = $false
$c.Characters(275,55).Font.Bold = $false
$c.Characters(275,55).Font.Italic = $false

$ws.Rows.Item(7).RowHeight = 90
